$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "where," paragraphs: drop the (now superfluous) grammar-check markers
#    that wrapped "where" (w:proofErr gramStart/gramEnd) while keeping the
#    paragraph formatting and wording identical.
#    A direct Find/Replace on "where" leaves the gramStart marker behind
#    because it sits exactly at the paragraph's first character boundary,
#    so we first nudge a throwaway character in front of it (which pushes
#    the marker off that boundary) and then replace across that boundary
#    in one shot, which drops both markers cleanly.
# ---------------------------------------------------------------------------
function Remove-WhereGrammarMarkers {
    param($doc)
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text -eq "where,`r") {
            $s = $para.Range.Start
            $doc.Range($s, $s).InsertBefore("Z")
            $doc.Content.Find.Execute("Zwhere,", $false, $false, $false, $false, $false, $true, 1, $false, "where,", 2) | Out-Null
        }
    }
}
Remove-WhereGrammarMarkers -doc $d

# ---------------------------------------------------------------------------
# 2) OH paragraph (annual operating hours breakdown): add a hanging indent
#    and abbreviate the hours/day, days/week, weeks/year phrasing.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*OH*Annual operating hours of the equipment*") {
        $para.LeftIndent = 36
        $para.FirstLineIndent = -36
        break
    }
}

$d.Content.Find.Execute(" hours per day, ", $false, $false, $false, $false, $false, $true, 1, $false,
    " hrs/day, " + [char]9 + "    ", 2) | Out-Null
$d.Content.Find.Execute(" days per week, ", $false, $false, $false, $false, $false, $true, 1, $false,
    " days/wk, ", 2) | Out-Null
$d.Content.Find.Execute(" weeks per year", $false, $false, $false, $false, $false, $true, 1, $false,
    " wks/yr", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Implementation cost references paragraph: remove the gramStart/gramEnd
#    markers around "in order to" by replacing straight across them.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("process in order to determine", $false, $false, $false, $false, $false, $true, 1, $false,
    "process in order to determine", 2) | Out-Null

Write-Host "Edit complete"
